$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 14
$ws1.Range("F6").Value = 535
$ws1.Range("F7").Value = 1626
$ws1.Range("F8").Value = 9
$ws1.Range("F9").Value = 19
$ws1.Range("F10").Value = 1426
$ws1.Range("F11").Value = 121
$ws1.Range("F12").Value = 28
$ws1.Range("F13").Value = 363
$ws1.Range("F14").Value = 247
$ws1.Range("F15").Value = 183
$ws1.Range("F17").Value = 15
$ws1.Range("F19").Value = 257
$ws1.Range("F20").Value = 141
$ws1.Range("F21").Value = 210
$ws1.Range("F22").Value = 196

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 14
$ws4.Range("F6").Value = 535
$ws4.Range("F7").Value = 1626
$ws4.Range("F9").Value = 9
$ws4.Range("F10").Value = 19
$ws4.Range("F11").Value = 1426
$ws4.Range("F12").Value = 121
$ws4.Range("F13").Value = 28
$ws4.Range("F14").Value = 363
$ws4.Range("F15").Value = 247
$ws4.Range("F16").Value = 183
$ws4.Range("F18").Value = 15
$ws4.Range("F20").Value = 257
$ws4.Range("F21").Value = 141
$ws4.Range("F22").Value = 210
$ws4.Range("F23").Value = 196
